$d = $word.ActiveDocument

# --- 1) Header table: semester/year line -------------------------------
# "2º Semestre 2024"  ->  "1º Semestre 2025"
$d.Content.Find.Execute("2º Semestre 2024", $false, $false, $false, $false, $false,
                         $true, 1, $false, "1º Semestre 2025", 2) | Out-Null

# --- 2) Header table: date placeholder line -----------------------------
# "20 / 08 / 2024"  ->  "____ / ___ / 2025"
$d.Content.Find.Execute("20 / 08 / 2024", $false, $false, $false, $false, $false,
                         $true, 1, $false, "____ / ___ / 2025", 2) | Out-Null

# --- 3) EXE 004 statement: re-type the sentence about the song ----------
# Content is unchanged; the source edit simply retypes the sentence so the
# three runs that made it up collapse into a single run.
$sentence = "Pe" + [char]0x00E7 + "a ao usu" + [char]0x00E1 + "rio para digitar a primeira linha de uma m" + [char]0x00FA + "sica qualquer  e dever" + [char]0x00E1 + " exibir o comprimento dela. Pe" + [char]0x00E7 + "a um n" + [char]0x00FA + "mero inicial e um n" + [char]0x00FA + "mero final e, em seguida, exiba apenas essa se" + [char]0x00E7 + [char]0x00E3 + "o da letra da " + [char]0x006D + [char]0x00FA + "sica (lembre-se de que o Python come" + [char]0x00E7 + "a a contar a partir de 0 e n" + [char]0x00E3 + "o de 1)."
$d.Content.Find.Execute($sentence, $false, $false, $false, $false, $false,
                         $true, 1, $false, $sentence, 2) | Out-Null
